$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.418.90"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.724.56"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'243.41"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4908"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").Value = "'0.2609"
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "'0.06195"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "1.722.13"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'0.06997"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'15.51"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "'4.577"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'0.6000"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'77.24"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "26.415.92"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'0.000007142"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Value = "'11.35"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "1.947.84"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "'4.476"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").Value = "'8.573"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "'5.147"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "'137.52"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Value = "'1.397"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "'107.15"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'1.699"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("D30").Value = "'3.944"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "'0.07955"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'0.04532"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.603"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'0.9963"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6262"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9156"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.394"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'1.950"
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'0.9999"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01481"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'99.99"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.347"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3841"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.704"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1158"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05360"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.717"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.15"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.239"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'50.92"
$ws.Range("E51").Value = "  -0.42%  "